# Apply the documentation edit to the "Additional Design Work" textbox on
# slide 1: split the sentence about the 3D-printed enclosures into three
# runs (inserting "were " / "designed using ") and let the textbox's
# autosize grow to fit the extra line of text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 53")   # "Additional Design Work" box
$tf = $sh.TextFrame
$tr = $tf.TextRange

$oldPhrase = ". 3D-printed enclosures for both the LIDAR and boards Autodesk Fusion and printed on a "
$newPhrase = ". 3D-printed enclosures for both the LIDAR and boards were designed using Autodesk Fusion and printed on a "

$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldPhrase) + 1

# Replace the old phrase with the new, longer phrase (still a single run
# at this point).
$target = $tr.Characters($startPos, $oldPhrase.Length)
$origSize = $target.Font.Size
$target.Text = $newPhrase

# Split the merged run into three runs at the correct word boundaries by
# nudging (re-asserting) the character formatting of the middle span
# ("designed using "), which forces PowerPoint to give it its own <a:r>.
$fullText2 = $tr.Text
$midStart = $fullText2.IndexOf("designed using ") + 1
$midLen = "designed using ".Length
$midRange = $tr.Characters($midStart, $midLen)
$midRange.Font.Size = $origSize

# The textbox is set to auto-fit; re-flowing the extra line grows its
# height to match the new text content. (914400 EMU/in, 72 pt/in.)
$targetHeightEmu = 7976671
$sh.Height = $targetHeightEmu / 914400 * 72
